# Weekly price-sheet update: a new record for 2022-05-26 is inserted as the
# new row 8, pushing the existing rows 8-18 down to rows 9-19 (dimension
# therefore grows from A1:R18 to A1:R19). Rows 1-7 are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 - this shifts old rows 8..18 down to 9..19
# (carrying their values/styles with them), exactly matching the diff.
$ws.Rows.Item(8).Insert()

# Populate the freshly inserted row 8 with the new market record.
$ws.Range("A8").Value = 5
$ws.Range("B8").Value = "Macroferia Regional de Talca"
$ws.Range("C8").Value = "Maule"
$ws.Range("D8").Value = 44707
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 100112040
$ws.Range("G8").Value = "Cilantro"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 150
$ws.Range("K8").Value = 9000
$ws.Range("L8").Value = 9000
$ws.Range("M8").Value = 9000
$ws.Range("N8").Value = '$/caja 36 atados'
$ws.Range("O8").Value = "Región Metropolitana"
$ws.Range("P8").Value = 250
$ws.Range("Q8").Value = 36
$ws.Range("R8").Value = "Hortaliza"
